$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Read")

$ws.Range("A2:A22").UnMerge()
